# feat: add 2022-Q1 data
#
# The workbook previously had a single "总计" (totals) sheet summarizing
# quarterly holdings. This change:
#   1. Turns that "总计" sheet into a new "2022-Q1" detail sheet (same
#      column layout as the other quarterly detail sheets: fund code,
#      fund name, fund size, total stock position, position ratio,
#      held market value, position rank).
#   2. Adds a brand-new "总计" sheet right after it, which is the old
#      totals table with a new first row for 2022-Q1 prepended (and the
#      running index column renumbered).

$wb = $excel.ActiveWorkbook

# A donor cell that already carries the "header / index column" style
# (bold font, thin border, centered) used throughout the workbook, so we
# can stamp new cells with the exact same style index instead of growing
# the style table.
$styleDonorSheet = $wb.Worksheets.Item("2021-Q4")
$styleDonor = $styleDonorSheet.Range("B1")

function Stamp-HeaderStyle($range) {
    $styleDonor.Copy()
    $range.PasteSpecial(-4122) # xlPasteFormats
    $wb.Application.CutCopyMode = $false
}

function Set-TextValue($range, [string]$value) {
    # Assigning a numeric-looking string to .Value auto-converts it to a
    # number; forcing a text number format first (and clearing it again
    # afterwards via the Normal style) keeps it stored as text, matching
    # the source data which keeps these figures as strings.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# Step 1: duplicate "总计" (so the fresh totals sheet keeps identical
# page setup / formatting to the original), then turn the original into
# the new "2022-Q1" detail sheet.
# ---------------------------------------------------------------------
$oldZj = $wb.Worksheets.Item("总计")
$oldZj.Copy($null, $oldZj) # places the duplicate right after $oldZj

$q1 = $oldZj
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"
Stamp-HeaderStyle $q1.Range("B1:H1")

$q1Rows = @(
    @("001985", "富国低碳新经济混合A",             "39.35", "93.88", "3.59", "1.4127", 7),
    @("519035", "富国天博创新混合",                 "25.63", "93.05", "3.16", "0.8099", 9),
    @("011357", "华泰柏瑞品质成长混合A",            "25.98", "76.67", "2.54", "0.6599", 4),
    @("006218", "富国生物医药科技混合A",            "9.55",  "87.25", "4.96", "0.4737", 4),
    @("009990", "华泰柏瑞品质优选混合A",            "13.12", "76.98", "2.76", "0.3621", 4),
    @("011466", "兴业医疗保健混合A",                 "7.10",  "84.68", "4.74", "0.3365", 4),
    @("000513", "富国高端制造行业股票",              "9.25",  "93.54", "3.62", "0.3348", 10),
    @("011921", "富国均衡成长三年持有期混合A",       "7.52",  "93.77", "3.18", "0.2391", 6),
    @("160529", "博时创业板两年定期开放混合",        "7.92",  "82.61", "2.64", "0.2091", 10),
    @("008138", "富国龙头优势混合",                   "5.02",  "93.93", "4.04", "0.2028", 6),
    @("100016", "富国天源沪港深平衡混合",             "6.23",  "72.29", "3.04", "0.1894", 5),
    @("011467", "兴业医疗保健混合C",                  "2.22",  "84.68", "4.74", "0.1052", 4),
    @("009991", "华泰柏瑞品质优选混合C",              "2.83",  "76.98", "2.76", "0.0781", 4),
    @("011308", "富国生物医药科技混合C",              "1.26",  "87.25", "4.96", "0.0625", 4),
    @("004558", "汇安丰裕灵活配置混合A",              "0.99",  "83.41", "4.35", "0.0431", 2),
    @("011358", "华泰柏瑞品质成长混合C",              "1.26",  "76.67", "2.54", "0.0320", 4),
    @("011922", "富国均衡成长三年持有期混合C",        "0.50",  "93.77", "3.18", "0.0159", 6),
    @("011306", "富国低碳新经济混合C",                "0.30",  "93.88", "3.59", "0.0108", 7),
    @("002802", "广发东财大数据精选灵活配置混合",     "0.41",  "55.13", "1.85", "0.0076", 9),
    @("004559", "汇安丰裕灵活配置混合C",              "0.01",  "83.41", "4.35", "0.0004", 2)
)

for ($i = 0; $i -lt $q1Rows.Count; $i++) {
    $r = $i + 2
    $data = $q1Rows[$i]

    $q1.Range("A$r").Value = $i
    Stamp-HeaderStyle $q1.Range("A$r")

    Set-TextValue $q1.Range("B$r") $data[0]
    $q1.Range("C$r").Value = $data[1]
    Set-TextValue $q1.Range("D$r") $data[2]
    Set-TextValue $q1.Range("E$r") $data[3]
    Set-TextValue $q1.Range("F$r") $data[4]
    Set-TextValue $q1.Range("G$r") $data[5]
    $q1.Range("H$r").Value = $data[6]
}

# ---------------------------------------------------------------------
# Step 2: the duplicated sheet becomes the new "总计" totals sheet,
# placed right after "2022-Q1".
# ---------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计 (2)")
$zj.Name = "总计"
$zj.Cells.Clear()

$zj.Range("B1").Value = "日期"
$zj.Range("C1").Value = "持有数量(只)"
$zj.Range("D1").Value = "持有市值(亿元)"
Stamp-HeaderStyle $zj.Range("B1:D1")

$zjRows = @(
    @("2022-Q1", 20, 5.59),
    @("2021-Q4", 33, 9.18),
    @("2021-Q3", 27, 6.59),
    @("2021-Q2", 18, 5.18),
    @("2021-Q1", 3,  0.33),
    @("2020-Q4", 2,  0.09)
)

for ($i = 0; $i -lt $zjRows.Count; $i++) {
    $r = $i + 2
    $data = $zjRows[$i]

    $zj.Range("A$r").Value = $i
    Stamp-HeaderStyle $zj.Range("A$r")

    $zj.Range("B$r").Value = $data[0]
    $zj.Range("C$r").Value = $data[1]
    $zj.Range("D$r").Value = $data[2]
}

# Restore the originally active sheet/tab (this edit doesn't change the
# workbook's active-sheet selection).
$wb.Worksheets.Item("2020-Q4").Activate()

